# Updates the cryptocurrency symbol list (Price / Volume(1h) columns) and
# re-orders the GateToken/KuCoinToken/MXToken/... block of rows, matching
# the upstream "Updated symbol list" data-refresh commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B (Coin) and C (Link) are plain, non-numeric-looking text, so a
# normal .Value assignment keeps them as text with no special handling.
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"

# Columns D (Price) and E (Volume(1h)) hold numeric-looking text such as
# "307.91" or "0.44%". A direct .Value assignment would make Excel
# auto-convert these into real numbers/percentages (and silently round
# or drop significant trailing zeros), which does not match the source
# workbook where these cells are stored as plain text. To keep them as
# text: force the Text number format, assign the value, then restore the
# Normal style so no stray number formatting is left on the cell.
$numericTextCells = @("D2","E2","D3","E3","E4","D5","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","D19","E19","E20","D21","E21","D22","E22","D23","E23","D24","E24","D25","E25","D27","E27","D39","E39","E40","D41","E41","D42","E42","D43","E43","D44","E44","D45","E45","E46","E47","D48","E48","E49","D50","E50","D51","E51")
foreach ($cellRef in $numericTextCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "307.91"
$ws.Range("E2").Value = "0.44%"
$ws.Range("D3").Value = "36.32"
$ws.Range("E3").Value = "0.95%"
$ws.Range("E4").Value = "0.56%"
$ws.Range("D5").Value = "0.08137"
$ws.Range("E5").Value = "0.62%"
$ws.Range("D6").Value = "2.078"
$ws.Range("E6").Value = "8.63%"
$ws.Range("D7").Value = "4.156"
$ws.Range("E7").Value = "0.21%"
$ws.Range("D8").Value = "7.869"
$ws.Range("E8").Value = "0.09%"
$ws.Range("D9").Value = "0.9291"
$ws.Range("E9").Value = "-0.20%"
$ws.Range("D10").Value = "0.1462"
$ws.Range("E10").Value = "16.59%"
$ws.Range("D11").Value = "0.1923"
$ws.Range("E11").Value = "0.35%"
$ws.Range("D12").Value = "0.09123"
$ws.Range("E12").Value = "-1.15%"
$ws.Range("D13").Value = "0.03457"
$ws.Range("E13").Value = "-1.45%"
$ws.Range("D14").Value = "0.09871"
$ws.Range("E14").Value = "-0.53%"
$ws.Range("D15").Value = "0.001405"
$ws.Range("E15").Value = "-1.32%"
$ws.Range("D16").Value = "0.006211"
$ws.Range("E16").Value = "-7.01%"
$ws.Range("D17").Value = "3.840"
$ws.Range("E17").Value = "6.40%"
$ws.Range("D18").Value = "3.400"
$ws.Range("E18").Value = "10.84%"
$ws.Range("D19").Value = "0.3463"
$ws.Range("E19").Value = "0.76%"
$ws.Range("E20").Value = "1.83%"
$ws.Range("D21").Value = "4.826"
$ws.Range("E21").Value = "-6.84%"
$ws.Range("D22").Value = "0.2337"
$ws.Range("E22").Value = "-7.69%"
$ws.Range("D23").Value = "0.04381"
$ws.Range("E23").Value = "-0.93%"
$ws.Range("D24").Value = "0.001234"
$ws.Range("E24").Value = "-0.06%"
$ws.Range("D25").Value = "0.004198"
$ws.Range("E25").Value = "-11.21%"
$ws.Range("D27").Value = "0.0001300"
$ws.Range("E27").Value = "-0.05%"
$ws.Range("D39").Value = "0.02047"
$ws.Range("E39").Value = "4.16%"
$ws.Range("E40").Value = "-1.86%"
$ws.Range("D41").Value = "0.007480"
$ws.Range("E41").Value = "-1.09%"
$ws.Range("D42").Value = "0.01013"
$ws.Range("E42").Value = "-0.18%"
$ws.Range("D43").Value = "0.1378"
$ws.Range("E43").Value = "0.43%"
$ws.Range("D44").Value = "0.002130"
$ws.Range("E44").Value = "1.38%"
$ws.Range("D45").Value = "0.009720"
$ws.Range("E45").Value = "-9.12%"
$ws.Range("E46").Value = "-1.18%"
$ws.Range("E47").Value = "-0.07%"
$ws.Range("D48").Value = "64.84"
$ws.Range("E48").Value = "-0.58%"
$ws.Range("E49").Value = "-3.78%"
$ws.Range("D50").Value = "0.00002099"
$ws.Range("E50").Value = "-0.07%"
$ws.Range("D51").Value = "0.0001999"
$ws.Range("E51").Value = "-0.07%"

foreach ($cellRef in $numericTextCells) {
    $ws.Range($cellRef).Style = "Normal"
}
